$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: clear the "true" value
$ws.Range("B7").Value = ""

# Date: 2025-06-28 -> 2025-11-18
# Force text entry (avoid Excel's date auto-detection), then restore the
# cell's original formatting (a neighboring cell in the same column keeps
# the untouched style) so no stray direct-format style gets stamped here.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2025-11-18"
$ws.Range("B8").ClearFormats()
$ws.Range("B9").Copy()
$ws.Range("B8").PasteSpecial(-4122)
